$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 13 (weekly refresh: a new week's price observation is
# prepended to the top of the data block, pushing the existing rows 13-49
# down to 14-50).
$ws.Rows("13:13").Insert()

# Populate the new row 13 with this week's data.
$ws.Range("A13").Value = 4
$ws.Range("B13").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C13").Value = "Los Lagos"
$ws.Range("D13").Value = 44497
$ws.Range("E13").Value = 10
$ws.Range("F13").Value = 100112026
$ws.Range("G13").Value = "Haba"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 70
$ws.Range("K13").Value = 10000
$ws.Range("L13").Value = 10000
$ws.Range("M13").Value = 10000
$ws.Range("N13").Value = "$/saco 25 kilos"
$ws.Range("O13").Value = "Región Metropolitana"
$ws.Range("P13").Value = 400
$ws.Range("Q13").Value = 25
$ws.Range("R13").Value = "Hortaliza"
